$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above row 398, shifting rows 398:440 down to 399:441
$ws.Rows.Item(398).Insert()

# Populate the newly inserted row 398 with the new price entry
$ws.Cells.Item(398, 1).Value = 10
$ws.Cells.Item(398, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(398, 3).Value = "La Araucanía"
$ws.Cells.Item(398, 4).Value = 44449
$ws.Cells.Item(398, 4).Style = $ws.Cells.Item(399, 4).Style
$ws.Cells.Item(398, 4).NumberFormat = $ws.Cells.Item(399, 4).NumberFormat
$ws.Cells.Item(398, 5).Value = 9
$ws.Cells.Item(398, 6).Value = 100112021
$ws.Cells.Item(398, 7).Value = "Ají"
$ws.Cells.Item(398, 8).Value = "Inferno"
$ws.Cells.Item(398, 9).Value = "Primera"
$ws.Cells.Item(398, 10).Value = 95
$ws.Cells.Item(398, 11).Value = 45000
$ws.Cells.Item(398, 12).Value = 45000
$ws.Cells.Item(398, 13).Value = 45000
$ws.Cells.Item(398, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(398, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(398, 16).Value = 3000
$ws.Cells.Item(398, 17).Value = 15
$ws.Cells.Item(398, 18).Value = "Hortaliza"
